# Add a "Dropdown column" backed by a hidden reference sheet listing cities.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the hidden reference sheet right after TEST_SHEET ---------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "List_reference_hidden_sheet"

# Populate the hidden sheet with the list of cities used by the dropdown.
$ws2.Range("A1").Value = "Delhi"
$ws2.Range("A2").Value = "Kolkata"
$ws2.Range("A3").Value = "Chennai"
$ws2.Range("A4").Value = "Asam"
$ws2.Range("A5").Value = "Udisha"
$ws2.Range("A6").Value = "Mumbai"
$ws2.Range("A7").Value = "Panjab"
$ws2.Range("A8").Value = "Shilong"

# Match the column widths used on the main sheet.
$ws2.Columns.Item(1).ColumnWidth = 30.42
$ws2.Columns.Item(2).ColumnWidth = 30.42

# --- Defined name pointing at the hidden sheet's list ---------------------
$wb.Names.Add("HiddenList", "=List_reference_hidden_sheet!`$A`$2:`$A`$9")

# --- Data validation dropdown on the main sheet, column G, rows 2-9 -------
$range = $ws1.Range("G2:G9")
$range.Validation.Add(3, 1, 1, "=HiddenList")
$range.Validation.IgnoreBlank = $true
$range.Validation.ShowInput = $false
$range.Validation.ShowError = $false

# Keep the original sheet as the active one.
$ws1.Activate()
